$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.943499255685424
$ws.Range("C2").Value = 0.2867880669805913
$ws.Range("D2").Value = 0.07942334282255104
$ws.Range("E2").Value = 0.422417126820946
$ws.Range("G2").Value = 0.002386298674390978
$ws.Range("I2").Value = 0.3794055995484555
$ws.Range("N2").Value = 0.8453349042530007
$ws.Range("O2").Value = 2.068135879610168

$ws.Range("B3").Value = 0.8302944320977872
$ws.Range("C3").Value = 0.2499844811424055
$ws.Range("D3").Value = 0.0719111095230005
$ws.Range("E3").Value = 0.3683573400674049
$ws.Range("G3").Value = 0.002389893946760261
$ws.Range("I3").Value = 0.3808875235087257
$ws.Range("N3").Value = 0.8542444459377236
$ws.Range("O3").Value = 2.028164200781674

$ws.Range("B4").Value = 0.7607646330926059
$ws.Range("C4").Value = 0.2273288960323328
$ws.Range("D4").Value = 0.06733570831745794
$ws.Range("E4").Value = 0.3352707887242445
$ws.Range("G4").Value = 0.002392216674147907
$ws.Range("I4").Value = 0.38214183791602
$ws.Range("N4").Value = 0.8601306423851582
$ws.Range("O4").Value = 2.005470027344842

$ws.Range("B5").Value = 0.7324258458006057
$ws.Range("C5").Value = 0.2180819254666062
$ws.Range("D5").Value = 0.06548045949212167
$ws.Range("E5").Value = 0.3218122999286948
$ws.Range("G5").Value = 0.002393192273751751
$ws.Range("I5").Value = 0.382739163897444
$ws.Range("N5").Value = 0.862633837071705
$ws.Range("O5").Value = 1.996683994677795

$ws.Range("B6").Value = 0.7277199380201296
$ws.Range("C6").Value = 0.2165455880056015
$ws.Range("D6").Value = 0.06517295421176073
$ws.Range("E6").Value = 0.3195789468982611
$ws.Range("G6").Value = 0.002393356030054544
$ws.Range("I6").Value = 0.3828435423489189
$ws.Range("N6").Value = 0.8630558033505409
$ws.Range("O6").Value = 1.995252903850599

$ws.Range("B7").Value = 0.7603824650507249
$ws.Range("C7").Value = 0.2272042474608043
$ws.Range("D7").Value = 0.06731065036932193
$ws.Range("E7").Value = 0.3350891863149599
$ws.Range("G7").Value = 0.002392229713567771
$ws.Range("I7").Value = 0.3821495452741281
$ws.Range("N7").Value = 0.860163978154425
$ws.Range("O7").Value = 2.005349668847401

$ws.Range("B8").Value = 0.9044710306848174
$ws.Range("C8").Value = 0.2741101978197378
$ws.Range("D8").Value = 0.07682536864560063
$ws.Range("E8").Value = 0.4037538904201483
$ws.Range("G8").Value = 0.002387514470269569
$ws.Range("I8").Value = 0.3798448303509332
$ws.Range("N8").Value = 0.848320627759712
$ws.Range("O8").Value = 2.053968249518419

$ws.Range("B9").Value = 1.186845533386077
$ws.Range("C9").Value = 0.3656387687408369
$ws.Range("D9").Value = 0.09578278595277823
$ws.Range("E9").Value = 0.5393549444235504
$ws.Range("G9").Value = 0.002379177610047134
$ws.Range("I9").Value = 0.3780774781122744
$ws.Range("N9").Value = 0.8283949931833305
$ws.Range("O9").Value = 2.164112463911408

$ws.Range("B10").Value = 1.394200985097768
$ws.Range("C10").Value = 0.4326240994646469
$ws.Range("D10").Value = 0.1099006663618809
$ws.Range("E10").Value = 0.63971586030992
$ws.Range("G10").Value = 0.002373600840036738
$ws.Range("I10").Value = 0.3784844421057372
$ws.Range("N10").Value = 0.8157683082808163
$ws.Range("O10").Value = 2.254259683327234

$ws.Range("B11").Value = 1.488512464733105
$ws.Range("C11").Value = 0.4630447767865462
$ws.Range("D11").Value = 0.1163662022844392
$ws.Range("E11").Value = 0.6855666833004506
$ws.Range("G11").Value = 0.002371181547430674
$ws.Range("I11").Value = 0.3790459360244682
$ws.Range("N11").Value = 0.8104616099283533
$ws.Range("O11").Value = 2.297317024447807

$ws.Range("B12").Value = 1.524223311657522
$ws.Range("C12").Value = 0.474557099041192
$ws.Range("D12").Value = 0.1188208587094408
$ws.Range("E12").Value = 0.7029600468829358
$ws.Range("G12").Value = 0.002370282233882238
$ws.Range("I12").Value = 0.3793131651248558
$ws.Range("N12").Value = 0.8085150361207809
$ws.Range("O12").Value = 2.313919724318964

$ws.Range("B13").Value = 1.516532478984004
$ws.Range("C13").Value = 0.4720780395673501
$ws.Range("D13").Value = 0.1182919231244739
$ws.Range("E13").Value = 0.6992126722645366
$ws.Range("G13").Value = 0.002370475170415554
$ws.Range("I13").Value = 0.3792531763237932
$ws.Range("N13").Value = 0.8089314641636207
$ws.Range("O13").Value = 2.310330739208553

$ws.Range("B14").Value = 1.491450478523461
$ws.Range("C14").Value = 0.4639920490412806
$ws.Range("D14").Value = 0.1165680216559792
$ws.Range("E14").Value = 0.6869970173434012
$ws.Range("G14").Value = 0.002371107223553089
$ws.Range("I14").Value = 0.3790668240096835
$ws.Range("N14").Value = 0.8103002018651466
$ws.Range("O14").Value = 2.298676951467542

$ws.Range("B15").Value = 1.476086637771232
$ws.Range("C15").Value = 0.4590381935356618
$ws.Range("D15").Value = 0.1155129048507035
$ws.Range("E15").Value = 0.6795186409699312
$ws.Range("G15").Value = 0.002371496562191922
$ws.Range("I15").Value = 0.3789598032162473
$ws.Range("N15").Value = 0.8111467947538955
$ws.Range("O15").Value = 2.291577549879889

$ws.Range("B16").Value = 1.388037141137602
$ws.Range("C16").Value = 0.4306350085184363
$ws.Range("D16").Value = 0.1094790053494421
$ws.Range("E16").Value = 0.6367235565078886
$ws.Range("G16").Value = 0.002373761305283963
$ws.Range("I16").Value = 0.3784553617629598
$ws.Range("N16").Value = 0.8161239188498755
$ws.Range("O16").Value = 2.251487286160852

$ws.Range("B17").Value = 1.334017191634018
$ws.Range("C17").Value = 0.413197470341288
$ws.Range("D17").Value = 0.1057885354082799
$ws.Range("E17").Value = 0.610522125666094
$ws.Range("G17").Value = 0.00237518070855364
$ws.Range("I17").Value = 0.3782426416177316
$ws.Range("N17").Value = 0.8192892682931472
$ws.Range("O17").Value = 2.227420244617036

$ws.Range("B18").Value = 1.302944917733782
$ws.Range("C18").Value = 0.4031630291570991
$ws.Range("D18").Value = 0.1036699399352727
$ws.Range("E18").Value = 0.5954700689153611
$ws.Range("G18").Value = 0.002376008187161415
$ws.Range("I18").Value = 0.378155686713967
$ws.Range("N18").Value = 0.8211510436481078
$ws.Range("O18").Value = 2.213770180714732

$ws.Range("B19").Value = 1.2924241538459
$ws.Range("C19").Value = 0.399764708836301
$ws.Range("D19").Value = 0.1029533156859799
$ws.Range("E19").Value = 0.5903767724050653
$ws.Range("G19").Value = 0.002376290262167801
$ws.Range("I19").Value = 0.3781323099280627
$ws.Range("N19").Value = 0.8217884739590033
$ws.Range("O19").Value = 2.209181496766888

$ws.Range("B20").Value = 1.339767859825201
$ws.Range("C20").Value = 0.4150542256508629
$ws.Range("D20").Value = 0.1061809709176345
$ws.Range("E20").Value = 0.6133093996837289
$ws.Range("G20").Value = 0.002375028465234593
$ws.Range("I20").Value = 0.3782616192032506
$ws.Range("N20").Value = 0.8189480520071086
$ws.Range("O20").Value = 2.229962259555862

$ws.Range("B21").Value = 1.498817752685113
$ws.Range("C21").Value = 0.4663672989432825
$ws.Range("D21").Value = 0.1170742017275899
$ws.Range("E21").Value = 0.6905841997561737
$ws.Range("G21").Value = 0.002370921118524186
$ws.Range("I21").Value = 0.3791200744120573
$ws.Range("N21").Value = 0.8098964610748993
$ws.Range("O21").Value = 2.302091843156234

$ws.Range("B22").Value = 1.602749081831689
$ws.Range("C22").Value = 0.4998606714731864
$ws.Range("D22").Value = 0.124230306908089
$ws.Range("E22").Value = 0.7412675709405647
$ws.Range("G22").Value = 0.002368334729279999
$ws.Range("I22").Value = 0.379999622218449
$ws.Range("N22").Value = 0.8043477322941612
$ws.Range("O22").Value = 2.350970112401228

$ws.Range("B23").Value = 1.547280731021544
$ws.Range("C23").Value = 0.4819885169873714
$ws.Range("D23").Value = 0.1204075696042537
$ws.Range("E23").Value = 0.7141996310265029
$ws.Range("G23").Value = 0.002369706197152155
$ws.Range("I23").Value = 0.3795008880305204
$ws.Range("N23").Value = 0.807275582725353
$ws.Range("O23").Value = 2.324722826484106

$ws.Range("B24").Value = 1.33716803092517
$ws.Range("C24").Value = 0.414214815763728
$ws.Range("D24").Value = 0.106003541160149
$ws.Range("E24").Value = 0.6120492375042232
$ws.Range("G24").Value = 0.002375097258721986
$ws.Range("I24").Value = 0.3782529293889425
$ws.Range("N24").Value = 0.8191021851485729
$ws.Range("O24").Value = 2.228812434225347

$ws.Range("B25").Value = 1.110474718942555
$ws.Range("C25").Value = 0.3409248239227622
$ws.Range("D25").Value = 0.09062150725742413
$ws.Range("E25").Value = 0.5025535839006778
$ws.Range("G25").Value = 0.002381336209591533
$ws.Range("I25").Value = 0.3782581086203578
$ws.Range("N25").Value = 0.8334320900740977
$ws.Range("O25").Value = 2.132710178905597
